# Applies the "D suite" changes:
#  - Update rows 33-35 Results column (PASS -> SKIP)
#  - Append two new test case rows (36, 37) to the "Test Cases" sheet
#  - Update the view's top-left cell / selection range to match the new extent

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Rows 33-35: Results column changes from PASS to SKIP
$ws.Range("E33").Value = "SKIP"
$ws.Range("E34").Value = "SKIP"
$ws.Range("E35").Value = "SKIP"

# New row 36 (description string is added to the shared-string table before the TCID)
$ws.Range("C36").Value = "Verify that POST tab count getting increased while appreciate post from Record view page"
$ws.Range("A36").Value = "PublishedAPostLikeCountTest"
$ws.Range("B36").Value = "TBD"
$ws.Range("D36").Value = "Y"
$ws.Range("E36").Value = "SKIP"

# New row 37
$ws.Range("A37").Value = "PublishedAPostTimeStampTest"
$ws.Range("B37").Value = "TBD"
$ws.Range("C37").Value = "Verify that Created Post displayed as per System date"
$ws.Range("D37").Value = "Y"
$ws.Range("E37").Value = "PASS"

# Copy formatting from the row above to keep styling consistent
$ws.Range("A35:E35").Copy() | Out-Null
$ws.Range("A36:E36").PasteSpecial(-4122) | Out-Null
$ws.Range("A37:E37").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Select() | Out-Null

# Update the view to reflect the new data extent
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D2:D37").Select() | Out-Null
